$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '26.134.34'
$ws.Range("E2").Value = '  +1.35%  '
# Row 3
$ws.Range("D3").Value = '1.770.59'
$ws.Range("E3").Value = '  +1.46%  '
# Row 4
$ws.Range("E4").Value = '  -0.07%  '
# Row 5
$ws.Range("E5").Value = '  +0.37%  '
# Row 6
$ws.Range("D6").Value = '''0.9999'
$ws.Range("E6").Value = '  -0.03%  '
# Row 7
$ws.Range("D7").Value = '''0.5236'
$ws.Range("E7").Value = '  +3.67%  '
# Row 8
$ws.Range("D8").Value = '''0.2771'
$ws.Range("E8").Value = '  +4.41%  '
# Row 9
$ws.Range("D9").Value = '''40.56'
$ws.Range("E9").Value = '  -3.24%  '
# Row 10
$ws.Range("E10").Value = '  +1.19%  '
# Row 11
$ws.Range("D11").Value = '1.770.98'
$ws.Range("E11").Value = '  +1.58%  '
# Row 12
$ws.Range("D12").Value = '''16.01'
$ws.Range("E12").Value = '  +3.79%  '
# Row 13
$ws.Range("D13").Value = '''0.07036'
$ws.Range("E13").Value = '  +1.67%  '
# Row 14
$ws.Range("D14").Value = '''0.6497'
$ws.Range("E14").Value = '  +8.42%  '
# Row 15
$ws.Range("D15").Value = '''4.538'
$ws.Range("E15").Value = '  +0.77%  '
# Row 16
$ws.Range("D16").Value = '''78.35'
$ws.Range("E16").Value = '  +1.74%  '
# Row 17
$ws.Range("D17").Value = '''1.000'
$ws.Range("E17").Value = '  -0.01%  '
# Row 18
$ws.Range("D18").Value = '''1.0000'
$ws.Range("E18").Value = '  -0.01%  '
# Row 19
$ws.Range("D19").Value = '26.143.63'
$ws.Range("E19").Value = '  +1.35%  '
# Row 20
$ws.Range("D20").Value = '''11.78'
$ws.Range("E20").Value = '  +1.62%  '
# Row 21
$ws.Range("D21").Value = '''0.000006783'
$ws.Range("E21").Value = '  -1.13%  '
# Row 22
$ws.Range("D22").Value = '1.999.67'
$ws.Range("E22").Value = '  +1.65%  '
# Row 23
$ws.Range("E23").Value = '  +1.09%  '
# Row 24
$ws.Range("D24").Value = '''8.458'
$ws.Range("E24").Value = '  +3.92%  '
# Row 25
$ws.Range("D25").Value = '''5.218'
$ws.Range("E25").Value = '  +0.16%  '
# Row 26
$ws.Range("D26").Value = '''138.12'
$ws.Range("E26").Value = '  +0.49%  '
# Row 27
$ws.Range("B27").Value = 'LidoDAOToken'
$ws.Range("C27").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D27").Value = '''1.880'
$ws.Range("E27").Value = '  +3.37%  '
# Row 28
$ws.Range("B28").Value = 'Toncoin'
$ws.Range("C28").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D28").Value = '''1.491'
$ws.Range("E28").Value = '  -1.81%  '
# Row 29
$ws.Range("D29").Value = '''15.21'
$ws.Range("E29").Value = '  +1.40%  '
# Row 30
$ws.Range("D30").Value = '''103.16'
$ws.Range("E30").Value = '  -0.46%  '
# Row 31
$ws.Range("D31").Value = '''0.08430'
$ws.Range("E31").Value = '  +3.97%  '
# Row 32
$ws.Range("D32").Value = '''3.746'
$ws.Range("E32").Value = '  -0.65%  '
# Row 33
$ws.Range("D33").Value = '''3.479'
$ws.Range("E33").Value = '  +0.41%  '
# Row 34
$ws.Range("D34").Value = '''0.04477'
$ws.Range("E34").Value = '  -0.72%  '
# Row 35
$ws.Range("D35").Value = '''2.660'
$ws.Range("E35").Value = '  +0.53%  '
# Row 36
$ws.Range("D36").Value = '''1.011'
$ws.Range("E36").Value = '  +2.74%  '
# Row 37
$ws.Range("D37").Value = '''0.6103'
$ws.Range("E37").Value = '  +0.57%  '
# Row 38
$ws.Range("D38").Value = '''2.770'
$ws.Range("E38").Value = '  +3.57%  '
# Row 39
$ws.Range("D39").Value = '''2.001'
$ws.Range("E39").Value = '  +4.21%  '
# Row 40
$ws.Range("E40").Value = '  +2.65%  '
# Row 41
$ws.Range("D41").Value = '''1.002'
$ws.Range("E41").Value = '  +0.17%  '
# Row 42
$ws.Range("D42").Value = '''103.13'
$ws.Range("E42").Value = '  +0.72%  '
# Row 43
$ws.Range("D43").Value = '''0.3884'
$ws.Range("E43").Value = '  +1.76%  '
# Row 44
$ws.Range("E44").Value = '  +2.63%  '
# Row 45
$ws.Range("D45").Value = '''4.948'
$ws.Range("E45").Value = '  -2.62%  '
# Row 46
$ws.Range("D46").Value = '''0.05515'
$ws.Range("E46").Value = '  +2.89%  '
# Row 47
$ws.Range("D47").Value = '''6.428'
$ws.Range("E47").Value = '  +8.42%  '
# Row 48
$ws.Range("D48").Value = '''0.1122'
$ws.Range("E48").Value = '  +1.25%  '
# Row 49
$ws.Range("D49").Value = '''30.38'
$ws.Range("E49").Value = '  +0.70%  '
# Row 50
$ws.Range("D50").Value = '''52.99'
$ws.Range("E50").Value = '  +0.72%  '
# Row 51
$ws.Range("E51").Value = '  +0.69%  '
